# Applies the "Fitted model and more measures" edit to
# 'Filtering Method' sheet: a new intermediate weighing/date column is
# inserted into the small table at rows 11-16 (pushing the old D:F data to
# E:G), and a brand-new weighing comparison table is added at rows 18-21.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Filtering Method")

# ---------------------------------------------------------------------
# 1) Rows 11-16: make room for a new column D by moving the existing
#    D:F data (formulas + values) one column to the right (E:G), then
#    populate the new D column with the extra weighing values.
#    We rebuild right-to-left so we never clobber a value before it is
#    copied onward.
# ---------------------------------------------------------------------

# Row 11 (headers): F->G, E->F, D->E, then new D11
$ws.Range("G11").Value = $ws.Range("F11").Value
$ws.Range("F11").Value = $ws.Range("E11").Value
$ws.Range("E11").Value = $ws.Range("D11").Value

# New D11: second date, same date style as C11 (copy number format only)
$ws.Range("C11").Copy() | Out-Null
$ws.Range("D11").PasteSpecial(-4122) | Out-Null
$ws.Range("D11").Value = 42331

# Row 12 (first data row) - plain formulas, no shared refs
$ws.Range("G12").Formula = "=F12/0.02"
$ws.Range("F12").Formula = "=E12-`$G`$3"
$ws.Range("E12").Formula = "=C12-B12"
$ws.Range("D12").Value = 2.806

# Row 13 (second data row)
$ws.Range("G13").Formula = "=F13/0.02"
$ws.Range("F13").Formula = "=E13-`$G`$3"
$ws.Range("E13").Formula = "=C13-B13"
$ws.Range("D13").Value = 2.8010000000000002

# Row 14 (third data row)
$ws.Range("G14").Formula = "=F14/0.02"
$ws.Range("F14").Formula = "=E14-`$G`$3"
$ws.Range("E14").Formula = "=C14-B14"
$ws.Range("D14").Value = 3.0030000000000001

# Row 15 (Average row)
$ws.Range("G15").Formula = "=AVERAGE(G12:G14)"
$ws.Range("F15").Formula = "=AVERAGE(F12:F14)"
$ws.Range("E15").Formula = "=AVERAGE(E12:E14)"
$ws.Range("D15").ClearContents() | Out-Null

# Row 16 (Stdev row)
$ws.Range("G16").Formula = "=STDEV(G12:G14)"
$ws.Range("F16").Formula = "=STDEV(F12:F14)"
$ws.Range("E16").Formula = "=STDEV(E12:E14)"
$ws.Range("D16").ClearContents() | Out-Null

# ---------------------------------------------------------------------
# 2) Row 3 (chart source cells) now point one column further right,
#    since the averages/stdevs moved from column F to column G.
# ---------------------------------------------------------------------
$ws.Range("L3").Formula = "=G15"
$ws.Range("M3").Formula = "=G16"

# ---------------------------------------------------------------------
# 3) New comparison table at rows 18-21 ("Initial Wt" vs "Added").
# ---------------------------------------------------------------------
$ws.Range("A10").Copy() | Out-Null
$ws.Range("A18").PasteSpecial(-4122) | Out-Null
$ws.Range("A18").Value = 42327
$ws.Range("B18").Value = "Initial Wt"

$ws.Range("A10").Copy() | Out-Null
$ws.Range("C18").PasteSpecial(-4122) | Out-Null
$ws.Range("C18").Value = 42331
$ws.Range("D18").Value = "Added"

$ws.Range("A19").Value = 1
$ws.Range("B19").Value = 2.8159999999999998
$ws.Range("C19").Value = 2.831
$ws.Range("D19").Formula = "=C19-B19"

$ws.Range("A20").Value = 2
$ws.Range("B20").Value = 3.2570000000000001
$ws.Range("C20").Value = 3.2709999999999999
$ws.Range("D20").Formula = "=C20-B20"

$ws.Range("A21").Value = 3
$ws.Range("B21").Value = 3.3109999999999999
$ws.Range("C21").Value = 3.3210000000000002
$ws.Range("D21").Formula = "=C21-B21"

Write-Output "Filtering Method sheet updated"
